$wb = $excel.ActiveWorkbook

# "Moving dates in UT to far future":
# The "Game" sheet stores the game start time as text in B3 ("Время старта").
# Push it 30 years out, from 2020 to 2050, keeping the same format/timezone.
$ws = $wb.Worksheets.Item("Game")
$ws.Range("B3").Value = "2050-07-01 00:00:00 +0300"

# Make "Game" the active/selected sheet (previously "Level 2" was active).
$ws.Activate()
